$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily TCM rows (06-08-2021 .. 06-09-2021) appended after existing row 150.
$data = @(
    @("06-08-2021", 113.64, 185.73, 100.4),
    @("09-08-2021", 114.44, 187.05, 101.04),
    @("10-08-2021", 114.78, 187.65, 101.32),
    @("11-08-2021", 113.7, 185.74, 100.37),
    @("12-08-2021", 112.78, 184.21, 99.56999999999999),
    @("13-08-2021", 112.7, 184.06, 99.5),
    @("16-08-2021", 113.07, 184.94, 99.84999999999999),
    @("17-08-2021", 114.47, 187.25, 101.08),
    @("18-08-2021", 114.94, 188.05, 101.43),
    @("19-08-2021", 114.71, 187.62, 101.24),
    @("20-08-2021", 114.61, 187.75, 101.06),
    @("23-08-2021", 114.08, 186.96, 100.59),
    @("24-08-2021", 113.91, 186.67, 100.5),
    @("25-08-2021", 114.05, 186.47, 100.68),
    @("26-08-2021", 114.25, 186.72, 100.88),
    @("27-08-2021", 114.34, 187.01, 100.91),
    @("30-08-2021", 114.59, 187.19, 101.21),
    @("31-08-2021", 113.98, 186.13, 100.67),
    @("01-09-2021", 113.43, 185.03, 100.22),
    @("02-09-2021", 112.33, 183.28, 99.26000000000001),
    @("03-09-2021", 112.52, 183.77, 99.43000000000001),
    @("06-09-2021", 112.39, 183.55, 99.34999999999999)
)

$startRow = 151
$endRow = $startRow + $data.Count - 1
$dateRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))

# Column A holds "dd-mm-yyyy" style text labels (same as the existing rows above).
# Pre-format as Text so Excel does not reinterpret them as date serials, then
# restore the default cell style so the new cells match the unstyled data rows.
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}

$dateRange.Style = "Normal"
